$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.069.49'
$ws.Range('D2').NumberFormat = "General"
$ws.Range('E2').Value = '  +5.41%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.923.77'
$ws.Range('D3').NumberFormat = "General"
$ws.Range('E3').Value = '  +2.65%  '

$ws.Range('E4').Value = '  -0.80%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '325.54'
$ws.Range('D5').NumberFormat = "General"
$ws.Range('E5').Value = '  +3.08%  '

$ws.Range('E6').Value = '  -0.73%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5166'
$ws.Range('D7').NumberFormat = "General"
$ws.Range('E7').Value = '  +1.63%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4003'
$ws.Range('D8').NumberFormat = "General"
$ws.Range('E8').Value = '  +2.65%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.08471'
$ws.Range('D9').NumberFormat = "General"
$ws.Range('E9').Value = '  +0.77%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '42.91'
$ws.Range('D10').NumberFormat = "General"
$ws.Range('E10').Value = '  +2.65%  '

$ws.Range('E11').Value = '  +1.68%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.336'
$ws.Range('D12').NumberFormat = "General"
$ws.Range('E12').Value = '  +1.90%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '21.19'
$ws.Range('D13').NumberFormat = "General"
$ws.Range('E13').Value = '  +3.95%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.918.91'
$ws.Range('D14').NumberFormat = "General"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.344'
$ws.Range('D15').NumberFormat = "General"
$ws.Range('E15').Value = '  +1.43%  '

$ws.Range('E16').Value = '  -0.80%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '94.36'
$ws.Range('D17').NumberFormat = "General"
$ws.Range('E17').Value = '  +3.42%  '

$ws.Range('E18').Value = '  +1.07%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06774'
$ws.Range('D19').NumberFormat = "General"
$ws.Range('E19').Value = '  +0.99%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '18.01'
$ws.Range('D20').NumberFormat = "General"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').NumberFormat = "General"
$ws.Range('E21').Value = '  -0.68%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.060'
$ws.Range('D22').NumberFormat = "General"
$ws.Range('E22').Value = '  +2.14%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '30.075.65'
$ws.Range('D23').NumberFormat = "General"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.19'
$ws.Range('D24').NumberFormat = "General"
$ws.Range('E24').Value = '  +1.00%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.201'
$ws.Range('D25').NumberFormat = "General"
$ws.Range('E25').Value = '  -1.46%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.141.51'
$ws.Range('D26').NumberFormat = "General"
$ws.Range('E26').Value = '  +2.77%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '160.07'
$ws.Range('D27').NumberFormat = "General"
$ws.Range('E27').Value = '  -0.92%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '20.99'
$ws.Range('D28').NumberFormat = "General"
$ws.Range('E28').Value = '  +1.74%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.466'
$ws.Range('D29').NumberFormat = "General"
$ws.Range('E29').Value = '  +4.81%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '129.14'
$ws.Range('D30').NumberFormat = "General"
$ws.Range('E30').Value = '  +2.44%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.078'
$ws.Range('D31').NumberFormat = "General"
$ws.Range('E31').Value = '  +3.39%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.1060'
$ws.Range('D32').NumberFormat = "General"
$ws.Range('E32').Value = '  +1.50%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.070'
$ws.Range('D33').NumberFormat = "General"
$ws.Range('E33').Value = '  +4.81%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.655'
$ws.Range('D34').NumberFormat = "General"
$ws.Range('E34').Value = '  +1.28%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.02502'
$ws.Range('D35').NumberFormat = "General"
$ws.Range('E35').Value = '  +1.99%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.06606'
$ws.Range('D36').NumberFormat = "General"
$ws.Range('E36').Value = '  +0.85%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.2226'
$ws.Range('D37').NumberFormat = "General"
$ws.Range('E37').Value = '  +3.04%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.245'
$ws.Range('D38').NumberFormat = "General"
$ws.Range('E38').Value = '  +4.62%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '9.011'
$ws.Range('D39').NumberFormat = "General"
$ws.Range('E39').Value = '  +1.73%  '

$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.201'
$ws.Range('D40').NumberFormat = "General"
$ws.Range('E40').Value = '  +2.54%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.6543'
$ws.Range('D41').NumberFormat = "General"
$ws.Range('E41').Value = '  +1.80%  '

$ws.Range('E42').Value = '  -0.77%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '11.40'
$ws.Range('D43').NumberFormat = "General"
$ws.Range('E43').Value = '  +2.48%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.6139'
$ws.Range('D44').NumberFormat = "General"
$ws.Range('E44').Value = '  +1.66%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.16'
$ws.Range('D45').NumberFormat = "General"
$ws.Range('E45').Value = '  +1.48%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.749'
$ws.Range('D46').NumberFormat = "General"
$ws.Range('E46').Value = '  +1.62%  '

$ws.Range('E47').Value = '  +2.30%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.243'
$ws.Range('D48').NumberFormat = "General"
$ws.Range('E48').Value = '  +2.20%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '125.59'
$ws.Range('D49').NumberFormat = "General"
$ws.Range('E49').Value = '  +2.89%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '79.47'
$ws.Range('D50').NumberFormat = "General"
$ws.Range('E50').Value = '  +3.63%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.145'
$ws.Range('D51').NumberFormat = "General"
$ws.Range('E51').Value = '  -2.63%  '
